# Update the two-digit multiplication problems to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @("66×92=", "23×40="),
    @("36×51=", "30×33="),
    @("18×40=", "98×83="),
    @("91×31=", "22×28="),
    @("53×16=", "91×91="),
    @("11×22=", "14×72="),
    @("49×27=", "38×62="),
    @("72×86=", "76×65="),
    @("59×57=", "87×19="),
    @("71×40=", "49×47="),
    @("20×26=", "72×88="),
    @("50×56=", "40×20="),
    @("54×71=", "26×30="),
    @("91×93=", "48×38="),
    @("96×89=", "15×14="),
    @("44×87=", "17×82="),
    @("38×32=", "89×64="),
    @("74×89=", "18×43="),
    @("55×73=", "41×13="),
    @("88×25=", "40×90="),
    @("20×90=", "97×69="),
    @("15×52=", "73×22="),
    @("37×66=", "43×81="),
    @("74×31=", "33×51="),
    @("97×99=", "36×48=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
